$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5373.2383
$ws.Range("J17").Value = 3641.9
$ws.Range("L17").Value = 10925.7
$ws.Range("N17").Value = -11261.7

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 35866.383
$ws.Range("I33").Value = 50177.867
$ws.Range("J33").Value = 87.666664
$ws.Range("K33").Value = 50177.867
$ws.Range("L33").Value = 87.666664
$ws.Range("M33").Value = -49948.867
$ws.Range("N33").Value = -545.666664

# ALC row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 4082
$ws.Range("J58").Value = 11500
$ws.Range("L58").Value = 34500
$ws.Range("N58").Value = -34800

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 50727.06
$ws.Range("J62").Value = 76662.71000000001
$ws.Range("L62").Value = 76662.71000000001
$ws.Range("N62").Value = -77910.71000000001

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 50727.06
$ws.Range("J65").Value = 76662.71000000001
$ws.Range("L65").Value = 383313.55
$ws.Range("N65").Value = -389553.55

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2093.0625
$ws.Range("J88").Value = 1982.4286
$ws.Range("L88").Value = 1982.4286
$ws.Range("N88").Value = -2794.4286

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 2093.0625
$ws.Range("J91").Value = 1982.4286
$ws.Range("L91").Value = 1982.4286
$ws.Range("N91").Value = -4790.4286

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 41764024
$ws.Range("J116").Value = 41690336
$ws.Range("L116").Value = 41690336
$ws.Range("N116").Value = -41697220

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1917.6818
$ws.Range("I138").Value = 1184.7693
$ws.Range("J138").Value = 2976.3333
$ws.Range("K138").Value = 3554.3079
$ws.Range("L138").Value = 8928.999899999999
$ws.Range("M138").Value = 1585.6921
$ws.Range("N138").Value = -19208.9999

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 5577.5
$ws.Range("I141").Value = 1656.7778
$ws.Range("K141").Value = 4970.3334
$ws.Range("M141").Value = 209.6665999999996

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2992.5
$ws.Range("I5").Value = 155.1
$ws.Range("J5").Value = 10086
$ws.Range("K5").Value = 155.1
$ws.Range("L5").Value = 10086
$ws.Range("M5").Value = -43.09999999999999
$ws.Range("N5").Value = -10310

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2613.9
$ws.Range("I74").Value = 2487.4119
$ws.Range("K74").Value = 2487.4119
$ws.Range("M74").Value = -1613.4119

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2613.9
$ws.Range("I77").Value = 2487.4119
$ws.Range("K77").Value = 12437.0595
$ws.Range("M77").Value = -8069.059499999999

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 8131925
$ws.Range("J122").Value = 1557.6
$ws.Range("L122").Value = 4672.799999999999
$ws.Range("N122").Value = -9572.799999999999

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2992.5
$ws.Range("I4").Value = 155.1
$ws.Range("J4").Value = 10086
$ws.Range("K4").Value = 155.1
$ws.Range("L4").Value = 10086
$ws.Range("M4").Value = -40.09999999999999
$ws.Range("N4").Value = -10316

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 77.7
$ws.Range("I7").Value = 76.75
$ws.Range("J7").Value = 81.5
$ws.Range("K7").Value = 76.75
$ws.Range("L7").Value = 81.5
$ws.Range("M7").Value = 36.25
$ws.Range("N7").Value = -307.5

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2196.7222
$ws.Range("I58").Value = 1163.7142
$ws.Range("K58").Value = 1163.7142
$ws.Range("M58").Value = -960.7141999999999

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5844
$ws.Range("I86").Value = 5876.769
$ws.Range("K86").Value = 5876.769
$ws.Range("M86").Value = -4753.769

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 5844
$ws.Range("I89").Value = 5876.769
$ws.Range("K89").Value = 29383.845
$ws.Range("M89").Value = -23767.845

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3621.7
$ws.Range("I132").Value = 3357.5557
$ws.Range("K132").Value = 10072.6671
$ws.Range("M132").Value = -7542.667099999999

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2196.7222
$ws.Range("I136").Value = 1163.7142
$ws.Range("K136").Value = 3491.1426
$ws.Range("M136").Value = -941.1425999999997

# CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 10227
$ws.Range("I87").Value = 666.6667
$ws.Range("J87").Value = 13413.777
$ws.Range("K87").Value = 2000.0001
$ws.Range("L87").Value = 40241.331
$ws.Range("M87").Value = -752.0001
$ws.Range("N87").Value = -42737.331

# CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 10227
$ws.Range("I90").Value = 666.6667
$ws.Range("J90").Value = 13413.777
$ws.Range("K90").Value = 6000.0003
$ws.Range("L90").Value = 120723.993
$ws.Range("M90").Value = 239.9997000000003
$ws.Range("N90").Value = -133203.993

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 129
$ws.Range("I2").Value = 8.75
$ws.Range("K2").Value = 8.75
$ws.Range("M2").Value = 104.25

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2733.5557
$ws.Range("I80").Value = 2816.6667
$ws.Range("J80").Value = 2692
$ws.Range("K80").Value = 2816.6667
$ws.Range("L80").Value = 2692
$ws.Range("M80").Value = -1818.6667
$ws.Range("N80").Value = -4688

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2733.5557
$ws.Range("I83").Value = 2816.6667
$ws.Range("J83").Value = 2692
$ws.Range("K83").Value = 14083.3335
$ws.Range("L83").Value = 13460
$ws.Range("M83").Value = -9091.333500000001
$ws.Range("N83").Value = -23444

# GSM row 138
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 75000
$ws.Range("J138").Value = 75000
$ws.Range("L138").Value = 75000
$ws.Range("N138").Value = -85280

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 618.5
$ws.Range("I16").Value = 585.7692
$ws.Range("J16").Value = 760.3333
$ws.Range("K16").Value = 585.7692
$ws.Range("L16").Value = 760.3333
$ws.Range("M16").Value = -415.7692
$ws.Range("N16").Value = -1100.3333

# WVR row 61
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 2500
$ws.Range("I61").Value = 2500
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2500
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2208
$ws.Range("N61").ClearContents()

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3392.3125
$ws.Range("I100").Value = 4106.077
$ws.Range("J100").Value = 299.33334
$ws.Range("K100").Value = 8212.154
$ws.Range("L100").Value = 598.66668
$ws.Range("M100").Value = -7671.154
$ws.Range("N100").Value = -1680.66668

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2443.147
$ws.Range("I136").Value = 988.9167
$ws.Range("K136").Value = 2966.7501
$ws.Range("M136").Value = -416.7501000000002
